$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 249.66667
$ws.Range("I12").Value = 274.5
$ws.Range("K12").Value = 274.5
$ws.Range("M12").Value = -104.5
$ws.Range("H87").Value = 38756.625
$ws.Range("J87").Value = 38756.625
$ws.Range("L87").Value = 38756.625
$ws.Range("N87").Value = -41252.625
$ws.Range("H90").Value = 38756.625
$ws.Range("J90").Value = 38756.625
$ws.Range("L90").Value = 116269.875
$ws.Range("N90").Value = -128749.875
$ws.Range("H98").Value = 2637.3
$ws.Range("I98").Value = 2789.8857
$ws.Range("J98").Value = 1569.2
$ws.Range("K98").Value = 2789.8857
$ws.Range("L98").Value = 1569.2
$ws.Range("M98").Value = -1291.8857
$ws.Range("N98").Value = -4565.2
$ws.Range("H107").Value = 2513.95
$ws.Range("I107").Value = 2642.2222
$ws.Range("J107").Value = 2409
$ws.Range("K107").Value = 2642.2222
$ws.Range("L107").Value = 2409
$ws.Range("M107").Value = -722.2222000000002
$ws.Range("N107").Value = -6249
$ws.Range("H122").Value = 2637.3
$ws.Range("I122").Value = 2789.8857
$ws.Range("J122").Value = 1569.2
$ws.Range("K122").Value = 8369.6571
$ws.Range("L122").Value = 4707.6
$ws.Range("M122").Value = -5919.6571
$ws.Range("N122").Value = -9607.6
$ws.Range("H131").Value = 266
$ws.Range("I131").Value = 266
$ws.Range("K131").Value = 798
$ws.Range("M131").Value = 4242
$ws.Range("H132").Value = 5957485.5
$ws.Range("I132").Value = 8551760
$ws.Range("J132").Value = 5914.647
$ws.Range("K132").Value = 25655280
$ws.Range("L132").Value = 17743.941
$ws.Range("M132").Value = -25652750
$ws.Range("N132").Value = -22803.941
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""
$ws.Range("H137").Value = 1181.1212
$ws.Range("I137").Value = 859.6923
$ws.Range("J137").Value = 1645.4073
$ws.Range("K137").Value = 2579.0769
$ws.Range("L137").Value = 4936.2219
$ws.Range("M137").Value = -29.07690000000002
$ws.Range("N137").Value = -10036.2219
$ws.Range("H138").Value = 1318.19
$ws.Range("I138").Value = 708.9394
$ws.Range("J138").Value = 1618.2687
$ws.Range("K138").Value = 2126.8182
$ws.Range("L138").Value = 4854.8061
$ws.Range("M138").Value = 3013.1818
$ws.Range("N138").Value = -15134.8061
$ws.Range("H141").Value = 719.3333
$ws.Range("I141").Value = 627.8570999999999
$ws.Range("K141").Value = 1883.5713
$ws.Range("M141").Value = 3296.4287

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5147.7144
$ws.Range("I32").Value = 4648.8447
$ws.Range("J32").Value = 7558.9165
$ws.Range("K32").Value = 4648.8447
$ws.Range("L32").Value = 7558.9165
$ws.Range("M32").Value = -4361.8447
$ws.Range("N32").Value = -8132.9165
$ws.Range("H110").Value = 1355.1111
$ws.Range("I110").Value = 927.3570999999999
$ws.Range("J110").Value = 1815.7693
$ws.Range("K110").Value = 927.3570999999999
$ws.Range("L110").Value = 1815.7693
$ws.Range("M110").Value = 1117.6429
$ws.Range("N110").Value = -5905.7693

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2286.652
$ws.Range("I31").Value = 2437.7222
$ws.Range("J31").Value = 1742.8
$ws.Range("K31").Value = 2437.7222
$ws.Range("L31").Value = 1742.8
$ws.Range("M31").Value = -2142.7222
$ws.Range("N31").Value = -2332.8
$ws.Range("H34").Value = 2286.652
$ws.Range("I34").Value = 2437.7222
$ws.Range("J34").Value = 1742.8
$ws.Range("K34").Value = 2437.7222
$ws.Range("L34").Value = 1742.8
$ws.Range("M34").Value = -2235.7222
$ws.Range("N34").Value = -2146.8
$ws.Range("H132").Value = 1895.3889
$ws.Range("I132").Value = 1541.52
$ws.Range("K132").Value = 4624.559999999999
$ws.Range("M132").Value = -2094.559999999999
$ws.Range("H135").Value = 30575
$ws.Range("J135").Value = 33514.285
$ws.Range("L135").Value = 33514.285
$ws.Range("N135").Value = -43654.285

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1116.7333
$ws.Range("I23").Value = 1933.3334
$ws.Range("J23").Value = 912.5833
$ws.Range("K23").Value = 5800.0002
$ws.Range("L23").Value = 2737.7499
$ws.Range("M23").Value = -5565.0002
$ws.Range("N23").Value = -3207.7499
$ws.Range("H61").Value = 324.2857
$ws.Range("I61").Value = 130
$ws.Range("J61").Value = 583.3333
$ws.Range("K61").Value = 390
$ws.Range("L61").Value = 1749.9999
$ws.Range("M61").Value = -175
$ws.Range("N61").Value = -2179.9999
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").Value = ""
$ws.Range("H131").Value = 27781596
$ws.Range("I131").Value = 166667170
$ws.Range("J131").Value = 4481.6333
$ws.Range("K131").Value = 500001510
$ws.Range("L131").Value = 13444.8999
$ws.Range("M131").Value = -499996470
$ws.Range("N131").Value = -23524.8999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7764706
$ws.Range("I11").Value = 7714286
$ws.Range("K11").Value = 7714286
$ws.Range("M11").Value = -7714147
$ws.Range("H21").Value = 1253950
$ws.Range("J21").Value = 5266.6665
$ws.Range("L21").Value = 5266.6665
$ws.Range("N21").Value = -5612.6665
$ws.Range("H30").Value = 1253950
$ws.Range("J30").Value = 5266.6665
$ws.Range("L30").Value = 5266.6665
$ws.Range("N30").Value = -5476.6665
$ws.Range("H80").Value = 2723.75
$ws.Range("I80").Value = 1767.7142
$ws.Range("K80").Value = 1767.7142
$ws.Range("M80").Value = -769.7141999999999
$ws.Range("H83").Value = 2723.75
$ws.Range("I83").Value = 1767.7142
$ws.Range("K83").Value = 8838.571
$ws.Range("M83").Value = -3846.571
$ws.Range("H97").Value = 621
$ws.Range("I97").Value = 644.8333
$ws.Range("J97").Value = 549.5
$ws.Range("K97").Value = 644.8333
$ws.Range("L97").Value = 549.5
$ws.Range("M97").Value = -148.8333
$ws.Range("N97").Value = -1541.5
$ws.Range("H113").Value = 1546.0834
$ws.Range("I113").Value = 1485.3
$ws.Range("K113").Value = 1485.3
$ws.Range("M113").Value = 684.7
$ws.Range("H132").Value = 2711.5
$ws.Range("I132").Value = 2271
$ws.Range("J132").Value = 3592.5
$ws.Range("K132").Value = 6813
$ws.Range("L132").Value = 10777.5
$ws.Range("M132").Value = -4283
$ws.Range("N132").Value = -15837.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1416.5264
$ws.Range("I7").Value = 1416.5264
$ws.Range("K7").Value = 1416.5264
$ws.Range("M7").Value = -1304.5264
$ws.Range("H22").Value = 599.8333
$ws.Range("I22").Value = 487.125
$ws.Range("J22").Value = 825.25
$ws.Range("K22").Value = 487.125
$ws.Range("L22").Value = 825.25
$ws.Range("M22").Value = -192.125
$ws.Range("N22").Value = -1415.25
$ws.Range("H27").Value = 599.8333
$ws.Range("I27").Value = 487.125
$ws.Range("J27").Value = 825.25
$ws.Range("K27").Value = 487.125
$ws.Range("L27").Value = 825.25
$ws.Range("M27").Value = -380.125
$ws.Range("N27").Value = -1039.25
$ws.Range("H46").Value = 6264.2856
$ws.Range("I46").Value = 425
$ws.Range("K46").Value = 425
$ws.Range("M46").Value = -237
$ws.Range("H55").Value = 388.13635
$ws.Range("I55").Value = 236.33333
$ws.Range("K55").Value = 236.33333
$ws.Range("M55").Value = -63.33332999999999
$ws.Range("H93").Value = 652.6667
$ws.Range("I93").Value = 563.2
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 563.2
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = 684.8
$ws.Range("N93").Value = -3596
$ws.Range("H100").Value = 1999.6666
$ws.Range("I100").Value = 1749.5
$ws.Range("K100").Value = 1749.5
$ws.Range("M100").Value = -1208.5
$ws.Range("H126").Value = 1416.5264
$ws.Range("I126").Value = 1416.5264
$ws.Range("K126").Value = 4249.5792
$ws.Range("M126").Value = -1779.5792
$ws.Range("H132").Value = 25103.627
$ws.Range("I132").Value = 1180.6522
$ws.Range("J132").Value = 52615.05
$ws.Range("K132").Value = 3541.9566
$ws.Range("L132").Value = 157845.15
$ws.Range("M132").Value = -1011.9566
$ws.Range("N132").Value = -162905.15

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H132").Value = 2229.121
$ws.Range("I132").Value = 1848.625
$ws.Range("J132").Value = 2587.2354
$ws.Range("K132").Value = 5545.875
$ws.Range("L132").Value = 7761.706200000001
$ws.Range("M132").Value = -3015.875
$ws.Range("N132").Value = -12821.7062
$ws.Range("H136").Value = 724.5925999999999
$ws.Range("I136").Value = 648.7619
$ws.Range("J136").Value = 990
$ws.Range("K136").Value = 2126.2857
$ws.Range("L136").Value = 2970
$ws.Range("M136").Value = 603.7143000000001
$ws.Range("N136").Value = -8070
